$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("info")
$ws.Range("A1").Value = ""
$ws.Range("B1").Value = ""
$ws.Range("C1").Value = ""
